$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 1.42
$ws.Range("D1").Value = 4.25
$ws.Range("E1").Value = 7.5

$ws.Range("C2").Value = 3.05
$ws.Range("D2").Value = 3.2

$ws.Range("C3").Value = 1.67
$ws.Range("D3").Value = 3.5
$ws.Range("E3").Value = 5

$ws.Range("C4").Value = 3.1
$ws.Range("D4").Value = 3.2
$ws.Range("E4").Value = 2.25

$ws.Range("C5").Value = 2.7
$ws.Range("D5").Value = 3.1
$ws.Range("E5").Value = 2.6

$ws.Range("C6").Value = 2.6
$ws.Range("D6").Value = 3.1
$ws.Range("E6").Value = 2.7

$ws.Range("C7").Value = 2.1
$ws.Range("D7").Value = 3.3
$ws.Range("E7").Value = 3.3

$ws.Range("C8").Value = 2.25
$ws.Range("D8").Value = 3.25
$ws.Range("E8").Value = 3.05

$ws.Range("C9").Value = 4.75
$ws.Range("D9").Value = 3.75
$ws.Range("E9").Value = 1.65

$ws.Range("C10").Value = 1.77
$ws.Range("D10").Value = 3.4
$ws.Range("E10").Value = 4.5

$ws.Range("E10").Select()
